$d = $word.ActiveDocument

# --- 1. Append a trailing space run to the end of the "What has been working?"
#        answer paragraph (right before its paragraph mark). ---
$pWorking = $d.Paragraphs(6)
$posTrail = $pWorking.Range.End - 1
$rTrail = $d.Range($posTrail, $posTrail)
$rTrail.InsertAfter(" ")

# --- 2. Fill in the empty answer paragraph under "What has not been working?"
#        This paragraph currently carries a bold/underlined "paragraph mark"
#        formatting (inherited from the heading above); replacing its XML via
#        InsertXML drops that mark and leaves a plain run, just like a user
#        turning off Bold/Underline before typing the answer. ---
$pNotWorking = $d.Paragraphs(9)
$rNotWorking = $pNotWorking.Range
$rNotWorking.Collapse(1)
$xmlNotWorking = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>We haven’t really thought about what trends we should focus on/look into.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rNotWorking.InsertXML($xmlNotWorking)

# --- 3. Fill in the empty answer paragraph under
#        "What can you do differently to address what hasn't been working?"
#        Three separate runs, same reasoning as above for the paragraph mark. ---
$pDiff = $d.Paragraphs(11)
$rDiff = $pDiff.Range
$rDiff.Collapse(1)
$xmlDiff = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">We should look at meteorological reports to get a better idea of what trends we should look for and try to </w:t></w:r><w:r><w:t>include those trends into our Trend class.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rDiff.InsertXML($xmlDiff)

# --- 4. Switch every paragraph's line spacing from 1.5 lines (360) to double
#        spacing (480) i.e. wdLineSpaceDouble, matching every <w:spacing> edit
#        in the diff. ---
foreach ($p in $d.Paragraphs) {
    $p.LineSpacingRule = 2
}
